$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Range("B9").Value = "09768edd95a8b219f10218dc50a94417"
$ws.Range("B11").Value = "b2b6ea8d6f2fd751653d2922bf86b7f7"
$ws.Range("B15").Value = "5b251fdfa08cb051878834729f44be74"
$ws.Range("B17").Value = "1ccfc1ec97dfed9f35c1ed5011b1cea9"
$ws.Range("B44").Value = "5d5ee50122e007df52ca745f7faf19bd"
$ws.Range("B89").Value = "f7945b435d376f43969ae850a7cc68cb"
$ws.Range("B99").Value = "45345d734b099da46e786c83e8f28c96"
$ws.Range("B110").Value = "74c498ae62afc36eaf69fb2be262b624"
$ws.Range("B121").Value = "16e942b2f0271e54d831782a253ff8bb"
$ws.Range("B133").Value = "0f5f13447ba864b1561c2ba55e4cd828"
$ws.Range("B136").Value = "eb7b0979e989c558249db2170fe6a48d"
$ws.Range("B154").Value = "6b15316edc1cc092b4abac42be90bd28"
$ws.Range("B159").Value = "4749c882ce4f82f5ec89fee91ecc415c"
$ws.Range("B160").Value = "a971ea9eb8c3823f3586968e3793190b"
$ws.Range("B168").Value = "92c63703e644491936dd6e9a8f2dc0c7"
$ws.Range("B169").Value = "4da83de0fa8baa0c3e34ef948fa497bf"
$ws.Range("B180").Value = "ae42a0af0e2092a422639ad4d71db265"
$ws.Range("B183").Value = "477b146f8b21754abe9e6418d07f97ae"
$ws.Range("B200").Value = "875decfdb4d3f6746c65a89f45459306"
$ws.Range("B222").Value = "b2c2d7b0c6e1e482e2baebfaa3e80238"
$ws.Range("B227").Value = "811e4b110a2cffba77fce045c7017d73"
$ws.Range("B228").Value = "5b813c348de89f8832b3df7554abeb70"
$ws.Range("B229").Value = "67e8de9238b1d980854c534789e8446c"
$ws.Range("B232").Value = "869c621bbced2dd1e9009bcaac137d49"
$ws.Range("B278").Value = "beba7bce29c4068483cd10898052ff4a"
$ws.Range("B335").Value = "7d3192fea74a6be1ead9e53c83c35f0f"
$ws.Range("B339").Value = "1eb832b6afed5fa4baf694d891211e50"
$ws.Range("B411").Value = "3168f1f5e791269381c5da9e31fafc82"
$ws.Range("B448").Value = "1566437fd9e351ae48076c533b2dc00e"
$ws.Range("B464").Value = "3c75af0a389448ba653dbb96b057f85d"
$ws.Range("B483").Value = "7c7e26fef28b133513b0e1d817db11ed"
$ws.Range("B507").Value = "bcf10a301975099317a3671d48f56727"
$ws.Range("B523").Value = "46abcc7d85f2732d753478da077c6dad"
$ws.Range("B542").Value = "caed40e30b8d326c9ee29159f49801d9"
$ws.Range("B561").Value = "6dae6fa19d878e3e786208dc34f13627"
$ws.Range("B580").Value = "90e9978e5fac4cdc1c413f6cc4049a3c"
$ws.Range("B592").Value = "0500c3294f2fe90971052abfee60871b"
$ws.Range("B624").Value = "a619418188285d32ee4afa2a1af3c1ad"
$ws.Range("B635").Value = "eff5797203762a41ac372a1640233c11"
$ws.Range("B637").Value = "dd6cf510f77b597f6665105ab8145b56"
$ws.Range("B657").Value = "fd75fa52fadd7dfc963bc94f149e1b82"
$ws.Range("B663").Value = "fe482945d81d149f47714c402a6d6be0"
$ws.Range("B708").Value = "c73244e4d02da93b2f5418460dd36c9d"
$ws.Range("B723").Value = "c5ee5882e46f01af84add9b219ddf0c2"
$ws.Range("B741").Value = "bb4d978bddca5a3b9e367b73036947b2"
$ws.Range("B776").Value = "0a647b4a3f32e50bca26867df944df5e"
$ws.Range("B819").Value = "d05f60cb7fe7ed68b218c83ac767a514"
$ws.Range("B823").Value = "828dfcdbe017b46b27ba6a91372baea2"
$ws.Range("B827").Value = "3cad1c31d6cda35f1ce8b17cbb9cfdb9"
$ws.Range("B833").Value = "55ee70e9919cf8142a528225a340560d"
$ws.Range("B835").Value = "e8dfad8ff97156163b1440cb8b6475c6"
$ws.Range("B838").Value = "e08d817cc6a46610a3b5f893585aa94e"
$ws.Range("B843").Value = "b102e7c044aa28ec0c96f4f071d794ab"
$ws.Range("B863").Value = "1b68267226727b46676d77be5d114c03"
